$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122
$xlPasteFormats = -4122

# Update existing row 228 (E, F, G values changed)
$ws.Range("E228").Value = 21.0002
$ws.Range("F228").Value = 21.35805
$ws.Range("G228").Value = 1474987

# Copy formatting of row 228 down to the three new rows so dates in
# column A keep the same date/time number format, font, border, alignment.
$ws.Range("A228:G228").Copy()
$ws.Range("A229:G231").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# New row 229
$ws.Range("A229").Value = 45047.29166666666
$ws.Range("B229").Value = "OANDA:USDCZK"
$ws.Range("C229").Value = 21.35805
$ws.Range("D229").Value = 22.32265
$ws.Range("E229").Value = 21.1644
$ws.Range("F229").Value = 22.1908
$ws.Range("G229").Value = 1568853

# New row 230
$ws.Range("A230").Value = 45078.29166666666
$ws.Range("B230").Value = "OANDA:USDCZK"
$ws.Range("C230").Value = 22.1908
$ws.Range("D230").Value = 22.24515
$ws.Range("E230").Value = 21.4818
$ws.Range("F230").Value = 21.7749
$ws.Range("G230").Value = 1549990

# New row 231
$ws.Range("A231").Value = 45110.29166666666
$ws.Range("B231").Value = "OANDA:USDCZK"
$ws.Range("C231").Value = 21.7749
$ws.Range("D231").Value = 22.0456
$ws.Range("E231").Value = 21.69175
$ws.Range("F231").Value = 21.73895
$ws.Range("G231").Value = 319645
